$d = $word.ActiveDocument

# 1. First card block: "{d[i].Last}, {d[i].First}  " -> "{d[i].LastFirstName}  "
$d.Content.Find.Execute("{d[i].Last}, {d[i].First}  {d[i].LicenceNumber}", $true, $false, $false, $false, $false, $true, 1, $false, "{d[i].LastFirstName}  {d[i].LicenceNumber}", 2)

# 2. "is duly licensed..." block for i: Start -> StartDate, End -> ExpiryDate
$d.Content.Find.Execute("for the period {d[i].Start} to {d[i].End}.", $true, $false, $false, $false, $false, $true, 1, $false, "for the period {d[i].StartDate} to {d[i].ExpiryDate}.", 2)

# 3. Second card block: "{d[i+1].Last}, {d[i+1].First}  " -> "{d[i+1].LastFirstName}  "
$d.Content.Find.Execute("{d[i+1].Last}, {d[i+1].First}  {d[i+1].LicenceNumber}", $true, $false, $false, $false, $false, $true, 1, $false, "{d[i+1].LastFirstName}  {d[i+1].LicenceNumber}", 2)

# 4. "...on the business of dealer agent for {d[i+1]..." block: Start -> StartDate (with trailing space), End -> ExpiryDate (with trailing space)
$d.Content.Find.Execute("for the period {d[i+1].Start} to {d[i+1].End}.", $true, $false, $false, $false, $false, $true, 1, $false, "for the period {d[i+1].StartDate } to {d[i+1].ExpiryDate }.", 2)
